# Apply update: all "Förändrad" dates in column C (rows 2-37) move from
# 45654 to 45655, and rows 36/37 swap their Beteckning (A) and Area (G)
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment column C (Förändrad) by 1 day for rows 2 through 37.
for ($r = 2; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}

# Rows 36 and 37 swap their Beteckning (A) and Area (G) values.
$a36 = $ws.Cells.Item(36, 1).Value2
$a37 = $ws.Cells.Item(37, 1).Value2
$ws.Cells.Item(36, 1).Value2 = $a37
$ws.Cells.Item(37, 1).Value2 = $a36

$g36 = $ws.Cells.Item(36, 7).Value2
$g37 = $ws.Cells.Item(37, 7).Value2
$ws.Cells.Item(36, 7).Value2 = $g37
$ws.Cells.Item(37, 7).Value2 = $g36
